$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin price/volume (and a symbol re-shuffle in rows 6-17)
# exactly as captured by the authoritative diff. A leading apostrophe is used
# so Excel keeps storing these numeric-looking values as literal text (matching
# the original inlineStr text cells) instead of coercing them to Number.
$ws.Cells.Item(2, 4).Value = "'309.73"
$ws.Cells.Item(2, 5).Value = "'0.75%"
$ws.Cells.Item(3, 4).Value = "'40.85"
$ws.Cells.Item(3, 5).Value = "'0.84%"
$ws.Cells.Item(4, 4).Value = "'5.124"
$ws.Cells.Item(4, 5).Value = "'1.31%"
$ws.Cells.Item(5, 4).Value = "'0.07656"
$ws.Cells.Item(5, 5).Value = "'0.90%"
$ws.Cells.Item(6, 2).Value = "'FTXToken"
$ws.Cells.Item(6, 3).Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Cells.Item(6, 4).Value = "'1.608"
$ws.Cells.Item(6, 5).Value = "'0.38%"
$ws.Cells.Item(7, 2).Value = "'MXToken"
$ws.Cells.Item(7, 3).Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(7, 4).Value = "'0.9119"
$ws.Cells.Item(7, 5).Value = "'1.08%"
$ws.Cells.Item(8, 4).Value = "'2.445"
$ws.Cells.Item(8, 5).Value = "'0.66%"
$ws.Cells.Item(9, 2).Value = "'LiechtensteinCryptoassetsExchange"
$ws.Cells.Item(9, 3).Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Cells.Item(9, 4).Value = "'0.1266"
$ws.Cells.Item(9, 5).Value = "'23.91%"
$ws.Cells.Item(10, 2).Value = "'WazirX"
$ws.Cells.Item(10, 3).Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Cells.Item(10, 4).Value = "'0.1807"
$ws.Cells.Item(10, 5).Value = "'2.94%"
$ws.Cells.Item(11, 2).Value = "'MandalaExchangeToken"
$ws.Cells.Item(11, 3).Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Cells.Item(11, 4).Value = "'0.09120"
$ws.Cells.Item(11, 5).Value = "'0.45%"
$ws.Cells.Item(12, 2).Value = "'BitrueCoin"
$ws.Cells.Item(12, 3).Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Cells.Item(12, 4).Value = "'0.04336"
$ws.Cells.Item(12, 5).Value = "'2.72%"
$ws.Cells.Item(13, 2).Value = "'BitMartToken"
$ws.Cells.Item(13, 3).Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Cells.Item(13, 4).Value = "'0.1046"
$ws.Cells.Item(13, 5).Value = "'-0.53%"
$ws.Cells.Item(14, 2).Value = "'BitForexToken"
$ws.Cells.Item(14, 3).Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Cells.Item(14, 4).Value = "'0.001246"
$ws.Cells.Item(14, 5).Value = "'1.39%"
$ws.Cells.Item(15, 2).Value = "'TigerCash"
$ws.Cells.Item(15, 3).Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Cells.Item(15, 4).Value = "'0.005775"
$ws.Cells.Item(15, 5).Value = "'-0.91%"
$ws.Cells.Item(16, 2).Value = "'LEO"
$ws.Cells.Item(16, 3).Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(16, 4).Value = "'3.352"
$ws.Cells.Item(16, 5).Value = "'0.04%"
$ws.Cells.Item(17, 2).Value = "'GateToken"
$ws.Cells.Item(17, 3).Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Cells.Item(17, 4).Value = "'4.281"
$ws.Cells.Item(17, 5).Value = "'0.18%"
$ws.Cells.Item(18, 5).Value = "'1.52%"
$ws.Cells.Item(19, 4).Value = "'6.914"
$ws.Cells.Item(19, 5).Value = "'2.26%"
$ws.Cells.Item(20, 5).Value = "'2.40%"
$ws.Cells.Item(21, 4).Value = "'0.2734"
$ws.Cells.Item(21, 5).Value = "'0.01%"
$ws.Cells.Item(22, 4).Value = "'0.04044"
$ws.Cells.Item(22, 5).Value = "'-3.56%"
$ws.Cells.Item(23, 4).Value = "'0.001269"
$ws.Cells.Item(23, 5).Value = "'3.09%"
$ws.Cells.Item(24, 4).Value = "'0.004043"
$ws.Cells.Item(24, 5).Value = "'-0.27%"
$ws.Cells.Item(25, 4).Value = "'0.0001269"
$ws.Cells.Item(25, 5).Value = "'-2.54%"
$ws.Cells.Item(26, 5).Value = "'24.20%"
$ws.Cells.Item(38, 4).Value = "'0.02424"
$ws.Cells.Item(38, 5).Value = "'2.14%"
$ws.Cells.Item(39, 4).Value = "'0.05254"
$ws.Cells.Item(39, 5).Value = "'2.00%"
$ws.Cells.Item(40, 4).Value = "'0.007839"
$ws.Cells.Item(40, 5).Value = "'1.21%"
$ws.Cells.Item(41, 4).Value = "'0.1303"
$ws.Cells.Item(41, 5).Value = "'1.40%"
$ws.Cells.Item(42, 4).Value = "'0.006797"
$ws.Cells.Item(42, 5).Value = "'-4.12%"
$ws.Cells.Item(43, 4).Value = "'0.001841"
$ws.Cells.Item(43, 5).Value = "'-4.26%"
$ws.Cells.Item(44, 4).Value = "'0.007436"
$ws.Cells.Item(44, 5).Value = "'-12.55%"
$ws.Cells.Item(45, 4).Value = "'0.3351"
$ws.Cells.Item(45, 5).Value = "'0.45%"
$ws.Cells.Item(46, 4).Value = "'0.00006873"
$ws.Cells.Item(46, 5).Value = "'7.73%"
$ws.Cells.Item(47, 4).Value = "'0.00000000750"
$ws.Cells.Item(47, 5).Value = "'-0.48%"
$ws.Cells.Item(48, 4).Value = "'0.1349"
$ws.Cells.Item(48, 5).Value = "'2,216.09%"
$ws.Cells.Item(49, 5).Value = "'-32.05%"
$ws.Cells.Item(50, 4).Value = "'0.00002099"
$ws.Cells.Item(50, 5).Value = "'-0.48%"
$ws.Cells.Item(51, 4).Value = "'0.0001999"
$ws.Cells.Item(51, 5).Value = "'-0.48%"
